$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these Price cells keep their literal text formatting (e.g. trailing zeros)
# (set individually - union ranges only apply NumberFormat to the first area)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = "56.704.09"
$ws.Range("E2").Value = "  +3.57%  "
$ws.Range("D3").Value = "3.247.36"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "395.34"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").Value = "108.95"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "0.587"
$ws.Range("E7").Value = "  +6.36%  "
$ws.Range("D8").Value = "3.245.21"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "0.627"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("D11").Value = "39.27"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "0.0981"
$ws.Range("E12").Value = "  +10.27%  "
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "3.762.45"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "8.36"
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("D16").Value = "19.12"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "3.252.37"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("E18").Value = "  -2.93%  "
$ws.Range("D19").Value = "10.76"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "56.647.53"
$ws.Range("E20").Value = "  +3.74%  "
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  +7.74%  "
$ws.Range("D23").Value = "12.99"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "295.60"
$ws.Range("E24").Value = "  +7.16%  "
$ws.Range("D25").Value = "74.42"
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("E26").Value = "  -2.89%  "
$ws.Range("D27").Value = "28.15"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").Value = "7.69"
$ws.Range("E29").Value = "  -4.31%  "
$ws.Range("D30").Value = "7.31"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "11.26"
$ws.Range("E34").Value = "  -3.14%  "
$ws.Range("D35").Value = "39.72"
$ws.Range("D36").Value = "0.0491"
$ws.Range("E36").Value = "  -2.92%  "
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("D38").Value = "51.33"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("D41").Value = "2.91"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "137.77"
$ws.Range("E42").Value = "  +5.17%  "
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "1.90"
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.99"
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("D46").Value = "17.10"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").Value = "0.282"
$ws.Range("E47").Value = "  -3.86%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("D50").Value = "2.151.14"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("E51").Value = "  -4.99%  "
